$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Tests" - insert two new columns (param:type, param:uuid) before
# the existing param:profile / param:lang columns, add their values, and add
# a new test row for "Missing Required Param".
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert two new columns at I:J (existing I/J -> K/L, keeping their widths).
$ws1.Range("I1:J1").EntireColumn.Insert()

# New header cells.
$ws1.Cells.Item(1, 9).Value = "param:type"
$ws1.Cells.Item(1, 10).Value = "param:uuid"

# New data cells for the existing test row (row 2).
$ws1.Cells.Item(2, 9).Value = "work"
$ws1.Cells.Item(2, 10).Value = "12345678-1234-1234-1234-123456789abc"

# New row 3: additional test case.
# NOTE: leading "'" forces these to be stored as text (matching the
# source workbook's t="str" cells) instead of being auto-typed as a
# boolean / left blank; resetting the Style afterwards drops the
# "quote prefix" marker so the cell format matches the rest of the sheet.
$ws1.Cells.Item(3, 1).Value = "get-data - Missing Required Param"
$ws1.Cells.Item(3, 2).Value = "Test GET /data/:type/:uuid with missing required parameters"
$ws1.Cells.Item(3, 3).Value = "'true"
$ws1.Cells.Item(3, 3).Style = "Normal"
$ws1.Cells.Item(3, 4).Value = 400
$ws1.Cells.Item(3, 5).Value = 10000
$ws1.Cells.Item(3, 6).Value = 2000
$ws1.Cells.Item(3, 7).Value = 500
$ws1.Cells.Item(3, 8).Value = "get-data,validation"
$ws1.Cells.Item(3, 9).Value = "'"
$ws1.Cells.Item(3, 9).Style = "Normal"
$ws1.Cells.Item(3, 10).Value = "12345678-1234-1234-1234-123456789abc"
$ws1.Cells.Item(3, 11).Value = "summary"
$ws1.Cells.Item(3, 12).Value = "en"

# ---------------------------------------------------------------------------
# Sheet 2: "Documentation" - document the new required path parameters.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Insert two new rows before the current param:profile row (row 18) for the
# new param:type / param:uuid descriptions.
$ws2.Range("A18:A19").EntireRow.Insert()

$ws2.Cells.Item(18, 1).Value = "param:type"
$ws2.Cells.Item(18, 2).Value = "type parameter (string) (REQUIRED - highlighted in yellow)"
$ws2.Cells.Item(19, 1).Value = "param:uuid"
$ws2.Cells.Item(19, 2).Value = "uuid parameter (string) (REQUIRED - highlighted in yellow)"

# Insert a new row after "Description:" note (row 26) for required parameters.
$ws2.Range("A27:A27").EntireRow.Insert()
$ws2.Cells.Item(27, 1).Value = "• Required parameters: type, uuid"
